# Update countries & provincias Spain
# Applies the data refresh that reorders a handful of countries within the
# ranking table (column A) and updates their associated statistics
# (columns B-H) to match the new day's figures. The sheet is sorted
# descending by "Casos totales" (column B), so a handful of countries
# changed rank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Rusia"
$ws.Range("B13").Value = 36793
$ws.Range("C13").Value = 4785
$ws.Range("D13").Value = 3057
$ws.Range("E13").Value = 33423
$ws.Range("F13").Value = 8
$ws.Range("G13").Value = 40
$ws.Range("H13").Value = 313

$ws.Range("A14").Value = "Belgica"
$ws.Range("B14").Value = 36138
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 7961
$ws.Range("E14").Value = 23014
$ws.Range("F14").Value = 1140
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 5163

$ws.Range("A15").Value = "Brasil"
$ws.Range("B15").Value = 34221
$ws.Range("C15").Value = 539
$ws.Range("D15").Value = 14026
$ws.Range("E15").Value = 18024
$ws.Range("F15").Value = 6634
$ws.Range("G15").Value = 30
$ws.Range("H15").Value = 2171

$ws.Range("A40").Value = "Singapur"
$ws.Range("B40").Value = 5992
$ws.Range("C40").Value = 942
$ws.Range("D40").Value = 708
$ws.Range("E40").Value = 5273
$ws.Range("F40").Value = 22
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 11

$ws.Range("A41").Value = "Indonesia"
$ws.Range("B41").Value = 5923
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 607
$ws.Range("E41").Value = 4796
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 520

$ws.Range("A42").Value = "Filipinas"
$ws.Range("B42").Value = 5878
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 487
$ws.Range("E42").Value = 5004
$ws.Range("F42").Value = 1
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 387

$ws.Range("A43").Value = "Serbia"
$ws.Range("B43").Value = 5690
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 534
$ws.Range("E43").Value = 5046
$ws.Range("F43").Value = 120
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 110

$ws.Range("A44").Value = "Malasia"
$ws.Range("B44").Value = 5251
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 2967
$ws.Range("E44").Value = 2198
$ws.Range("F44").Value = 51
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 86

$ws.Range("A45").Value = "Ucrania"
$ws.Range("B45").Value = 5106
$ws.Range("C45").Value = 444
$ws.Range("D45").Value = 275
$ws.Range("E45").Value = 4698
$ws.Range("F45").Value = 45
$ws.Range("G45").Value = 8
$ws.Range("H45").Value = 133

$ws.Range("A77").Value = "Oman"
$ws.Range("B77").Value = 1180
$ws.Range("C77").Value = 111
$ws.Range("D77").Value = 176
$ws.Range("E77").Value = 998
$ws.Range("F77").Value = 3
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 6

$ws.Range("A78").Value = "Republica de Macedonia"
$ws.Range("B78").Value = 1117
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 139
$ws.Range("E78").Value = 929
$ws.Range("F78").Value = 15
$ws.Range("G78").Value = 0
$ws.Range("H78").Value = 49

$ws.Range("A82").Value = "Afganistan"
$ws.Range("B82").Value = 933
$ws.Range("C82").Value = 27
$ws.Range("D82").Value = 112
$ws.Range("E82").Value = 791
$ws.Range("F82").Value = 0
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 30

$ws.Range("A83").Value = "Cuba"
$ws.Range("B83").Value = 923
$ws.Range("C83").Value = 0
$ws.Range("D83").Value = 192
$ws.Range("E83").Value = 700
$ws.Range("F83").Value = 17
$ws.Range("G83").Value = 0
$ws.Range("H83").Value = 31

